$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: "??" -> "Project Presentations" (B29) plus new "7 minute..." note (C29)
$ws.Range("B29").Value = "Project Presentations"
$ws.Range("C29").Value = "7 minute presentations + 7 minutes of questions per group"
$ws.Range("B26").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null

# Row 30: "??" -> "Course Summary Class" (B30) plus new cheatsheet link (D30)
$ws.Range("B30").Value = "Course Summary Class"
$ws.Range("D30").Value = "'- ``Summary Cheatsheet <https://github.com/nickeubank/unifyingdatascience/blob/master/lecture_slides/80_final_cheatsheet/UDS_Final_Cheatsheet.pdf>``_"
$ws.Range("D23").Copy() | Out-Null
$ws.Range("D30").PasteSpecial(-4122) | Out-Null

# Row 30 grows taller to fit the new wrapped text
$ws.Rows.Item(30).RowHeight = 68

# Move the active selection to B30 to match the saved view state
$ws.Range("B30").Select() | Out-Null
